# [ENH] budget report unit, add column
#
# The "Activity Group" column (P) is split into two columns:
#   - P10: "AG Code"  (renamed from "Activity Group")
#   - Q10: "AG Name"  (brand new column, inserted right after AG Code;
#                      everything that used to start at Q shifts one
#                      column to the right)
# A brand new "State" column is appended at the end of the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column Q (17th column), shifting
# "Job Order Code" ... "Planned Amount" (previously Q..AL) one column
# to the right (now R..AM). Formatting of the surrounding cells is
# carried along automatically by Insert().
$null = $ws.Range("Q1").EntireColumn.Insert()

# Rename the old "Activity Group" header and populate the new column.
$ws.Range("P10").Value = "AG Code"
$ws.Range("Q10").Value = "AG Name"

# Append the new "State" header at the end of row 10, copying the
# formatting (style) from the preceding header cell (AM10, "Planned
# Amount") so it matches the rest of the header row.
$null = $ws.Range("AM10").Copy()
$null = $ws.Range("AN10").PasteSpecial(-4122)
$ws.Range("AN10").Value = "State"
$excel.CutCopyMode = $false

# Match the saved selection/active cell from the edited workbook.
$null = $ws.Range("A9").Select()
